# Weekly fruit/vegetable price update:
# Insert a new observation row at row 196 (a new weekly price report for
# "Ají" variety "Inferno" in "Región de Arica y Parinacota"), pushing the
# existing rows 196:237 down to 197:238.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing rows 196:237 down by one to make room for the new record.
$ws.Rows(196).Insert()

# Populate the newly inserted row 196 with the new weekly report.
$ws.Range("A196").Value = 11
$ws.Range("B196").Value = "Vega Monumental Concepción"
$ws.Range("C196").Value = "Bíobío"
$ws.Range("D196").Value = 45218
$ws.Range("E196").Value = 8
$ws.Range("F196").Value = 100112021
$ws.Range("G196").Value = "Ají"
$ws.Range("H196").Value = "Inferno"
$ws.Range("I196").Value = "Primera"
$ws.Range("J196").Value = 30
$ws.Range("K196").Value = 38000
$ws.Range("L196").Value = 38000
$ws.Range("M196").Value = 38000
$ws.Range("N196").Value = "`$/caja 10 kilos"
$ws.Range("O196").Value = "Región de Arica y Parinacota"
$ws.Range("P196").Value = 3800
$ws.Range("Q196").Value = 10
$ws.Range("R196").Value = "Hortaliza"
